$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1134179.2
$ws.Range("J17").Value = 1231266
$ws.Range("L17").Value = 3693798
$ws.Range("N17").Value = -3694134

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 332.75
$ws.Range("I53").Value = 298.45456
$ws.Range("J53").Value = 374.66666
$ws.Range("K53").Value = 298.45456
$ws.Range("L53").Value = 374.66666
$ws.Range("M53").Value = 338.54544
$ws.Range("N53").Value = -1648.66666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 3000
$ws.Range("I54").Value = 3000
$ws.Range("K54").Value = 3000
$ws.Range("M54").Value = -2514

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 161.2
$ws.Range("I55").Value = 173.75
$ws.Range("J55").Value = 111
$ws.Range("K55").Value = 173.75
$ws.Range("L55").Value = 111
$ws.Range("M55").Value = 40.25
$ws.Range("N55").Value = -539

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 17867.334
$ws.Range("I94").Value = 17867.334
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 17867.334
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -17416.334
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 11212388
$ws.Range("I125").Value = 558.75
$ws.Range("J125").Value = 18686942
$ws.Range("K125").Value = 5028.75
$ws.Range("L125").Value = 168182478
$ws.Range("M125").Value = -2568.75
$ws.Range("N125").Value = -168187398

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 341120.6
$ws.Range("I132").Value = 380479.2
$ws.Range("J132").Value = 26251.75
$ws.Range("K132").Value = 1141437.6
$ws.Range("L132").Value = 78755.25
$ws.Range("M132").Value = -1138907.6
$ws.Range("N132").Value = -83815.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1662.3478
$ws.Range("I135").Value = 1715.6
$ws.Range("J135").Value = 1562.5
$ws.Range("K135").Value = 15440.4
$ws.Range("L135").Value = 14062.5
$ws.Range("M135").Value = -12905.4
$ws.Range("N135").Value = -19132.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 10000
$ws.Range("J43").Value = 10000
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10626

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 7400
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2464
$ws.Range("I61").Value = 1785.7778
$ws.Range("J61").Value = 4295.2
$ws.Range("K61").Value = 1785.7778
$ws.Range("L61").Value = 4295.2
$ws.Range("M61").Value = -1573.7778
$ws.Range("N61").Value = -4719.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1931.25
$ws.Range("I110").Value = 783.3333
$ws.Range("J110").Value = 2620
$ws.Range("K110").Value = 783.3333
$ws.Range("L110").Value = 2620
$ws.Range("M110").Value = 1261.6667
$ws.Range("N110").Value = -6710

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2215.1667
$ws.Range("I132").Value = 1930.8928
$ws.Range("J132").Value = 3210.125
$ws.Range("K132").Value = 5792.678400000001
$ws.Range("L132").Value = 9630.375
$ws.Range("M132").Value = -3262.678400000001
$ws.Range("N132").Value = -14690.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2464
$ws.Range("I136").Value = 1785.7778
$ws.Range("J136").Value = 4295.2
$ws.Range("K136").Value = 5357.3334
$ws.Range("L136").Value = 12885.6
$ws.Range("M136").Value = -2807.3334
$ws.Range("N136").Value = -17985.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4631.8335
$ws.Range("I31").Value = 1137.7407
$ws.Range("J31").Value = 15114.111
$ws.Range("K31").Value = 1137.7407
$ws.Range("L31").Value = 15114.111
$ws.Range("M31").Value = -842.7407000000001
$ws.Range("N31").Value = -15704.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4631.8335
$ws.Range("I34").Value = 1137.7407
$ws.Range("J34").Value = 15114.111
$ws.Range("K34").Value = 1137.7407
$ws.Range("L34").Value = 15114.111
$ws.Range("M34").Value = -935.7407000000001
$ws.Range("N34").Value = -15518.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 842.96
$ws.Range("I105").Value = 833.7
$ws.Range("J105").Value = 880
$ws.Range("K105").Value = 833.7
$ws.Range("L105").Value = 880
$ws.Range("M105").Value = 913.3
$ws.Range("N105").Value = -4374

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3584.16
$ws.Range("I132").Value = 3204.1177
$ws.Range("J132").Value = 4391.75
$ws.Range("K132").Value = 9612.3531
$ws.Range("L132").Value = 13175.25
$ws.Range("M132").Value = -7082.3531
$ws.Range("N132").Value = -18235.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5259
$ws.Range("I134").Value = 2209
$ws.Range("J134").Value = 6402.75
$ws.Range("K134").Value = 6627
$ws.Range("L134").Value = 19208.25
$ws.Range("M134").Value = -4092
$ws.Range("N134").Value = -24278.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4064.5
$ws.Range("J64").Value = 5000.6665
$ws.Range("L64").Value = 15001.9995
$ws.Range("N64").Value = -15541.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 4064.5
$ws.Range("J67").Value = 5000.6665
$ws.Range("L67").Value = 15001.9995
$ws.Range("N67").Value = -16873.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 21916.666
$ws.Range("J53").Value = 21916.666
$ws.Range("L53").Value = 21916.666
$ws.Range("N53").Value = -23178.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 222.125
$ws.Range("I107").Value = 175.28572
$ws.Range("J107").Value = 550
$ws.Range("K107").Value = 175.28572
$ws.Range("L107").Value = 550
$ws.Range("M107").Value = 1744.71428
$ws.Range("N107").Value = -4390

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2200
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2200
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6540

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2213.074
$ws.Range("I126").Value = 1713.1666
$ws.Range("J126").Value = 2355.9048
$ws.Range("K126").Value = 5139.4998
$ws.Range("L126").Value = 7067.714399999999
$ws.Range("M126").Value = -2669.4998
$ws.Range("N126").Value = -12007.7144

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2881.7073
$ws.Range("I132").Value = 2618.5806
$ws.Range("K132").Value = 7855.7418
$ws.Range("M132").Value = -5325.7418

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 6400
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 6400
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 6400
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -7436

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 7000
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 474.74075
$ws.Range("I55").Value = 353.2143
$ws.Range("J55").Value = 605.61536
$ws.Range("K55").Value = 353.2143
$ws.Range("L55").Value = 605.61536
$ws.Range("M55").Value = -180.2143
$ws.Range("N55").Value = -951.61536

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 22574.9
$ws.Range("J106").Value = 22574.9
$ws.Range("L106").Value = 22574.9
$ws.Range("N106").Value = -25098.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 13485.714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 31500
$ws.Range("J80").Value = 31500
$ws.Range("L80").Value = 31500
$ws.Range("N80").Value = -33496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 31500
$ws.Range("J83").Value = 31500
$ws.Range("L83").Value = 94500
$ws.Range("N83").Value = -104484

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 32818.184
$ws.Range("J123").Value = 32818.184
$ws.Range("L123").Value = 32818.184
$ws.Range("N123").Value = -42618.184
